$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 299, shifting existing rows 299:317 down to 300:318.
$ws.Range("A299:R299").EntireRow.Insert()

# Populate the newly inserted row 299 with the new data record.
$ws.Cells.Item(299, 1).Value = 3
$ws.Cells.Item(299, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(299, 3).Value = "Coquimbo"
$ws.Cells.Item(299, 4).Value = 44610
$ws.Cells.Item(299, 4).NumberFormat = $ws.Cells.Item(300, 4).NumberFormat
$ws.Cells.Item(299, 5).Value = 5
$ws.Cells.Item(299, 6).Value = 100112031
$ws.Cells.Item(299, 7).Value = "Poroto verde"
$ws.Cells.Item(299, 8).Value = "Magnum"
$ws.Cells.Item(299, 9).Value = "Primera"
$ws.Cells.Item(299, 10).Value = 35
$ws.Cells.Item(299, 11).Value = 25000
$ws.Cells.Item(299, 12).Value = 25000
$ws.Cells.Item(299, 13).Value = 25000
$ws.Cells.Item(299, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(299, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(299, 16).Value = 1000
$ws.Cells.Item(299, 17).Value = 25
$ws.Cells.Item(299, 18).Value = "Hortaliza"
